$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, shifting existing rows 233:320 down to 234:321
$ws.Rows("233:233").Insert()

# Populate the newly inserted row 233 with the new record's data
$ws.Cells.Item(233, 1).Value = 9
$ws.Cells.Item(233, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(233, 3).Value = "Metropolitana"
$ws.Cells.Item(233, 4).Value = 44917
$ws.Cells.Item(233, 5).Value = 13
$ws.Cells.Item(233, 6).Value = 100112001
$ws.Cells.Item(233, 7).Value = "Berenjena"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 90
$ws.Cells.Item(233, 11).Value = 14000
$ws.Cells.Item(233, 12).Value = 16000
$ws.Cells.Item(233, 13).Value = 15000
$ws.Cells.Item(233, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(233, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(233, 16).Value = 300
$ws.Cells.Item(233, 17).Value = 50
$ws.Cells.Item(233, 18).Value = "Hortaliza"
